# Update gh-pages output data (苏州-漫展信息.xlsx) — refresh "想去人数" (F column)
# counts on the "展览" (sheet 1) and "全部类型" (sheet 4) worksheets.

$wb = $excel.ActiveWorkbook

# sheet name -> { row -> new F value }
$sheetUpdates = @{
    "展览"   = @{
        2  = 609
        6  = 14255
        7  = 16201
        9  = 82
        21 = 133
        24 = 6493
        25 = 968
        26 = 11
        27 = 1112
        29 = 5685
        30 = 96
        32 = 168
        33 = 4715
    }
    "全部类型" = @{
        2  = 609
        6  = 14255
        7  = 16201
        9  = 82
        21 = 133
        25 = 6493
        26 = 968
        27 = 11
        28 = 1112
        31 = 5685
        32 = 96
        34 = 168
        35 = 4715
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetUpdates[$sheetName]
    foreach ($r in $rows.Keys) {
        $ws.Cells.Item($r, 6).Value = $rows[$r]
    }
}
